$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.819.41'
$ws.Cells.Item(2, 5).Value = '  -1.32%  '
$ws.Cells.Item(3, 4).Value = '1.634.94'
$ws.Cells.Item(3, 5).Value = '  -1.46%  '
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  -0.14%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '215.13'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.46%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '0.5010'
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -2.97%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -0.16%  '
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '0.2566'
$c.Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -0.89%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.06412'
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -0.53%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '19.55'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -1.96%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '0.07645'
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -1.81%  '
$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12, 4).Value = '1.636.73'
$ws.Cells.Item(12, 5).Value = '  -1.32%  '
$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '4.232'
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -1.51%  '
$ws.Cells.Item(14, 4).Value = '1.861.45'
$ws.Cells.Item(14, 5).Value = '  -1.40%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '0.5451'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -1.94%  '
$ws.Cells.Item(16, 4).Value = '0.0₅7911'
$ws.Cells.Item(16, 5).Value = '  -1.99%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '63.46'
$c.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -1.32%  '
$ws.Cells.Item(18, 4).Value = '25.863.45'
$ws.Cells.Item(18, 5).Value = '  -1.31%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -0.15%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '203.10'
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -3.98%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '4.301'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -2.70%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '9.942'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -1.00%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '5.981'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.23%  '
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '1.005'
$c.Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.10%  '
$ws.Cells.Item(25, 5).Value = '  +10.20%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '141.26'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -2.12%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '0.1145'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -1.55%  '
$ws.Cells.Item(28, 5).Value = '  -0.97%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '6.698'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -3.96%  '
$ws.Cells.Item(30, 5).Value = '  -1.33%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '0.04965'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -5.50%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '3.266'
$c.Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -2.96%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '3.181'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -1.32%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '1.529'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -2.62%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '2.353'
$c.Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -0.73%  '
$ws.Cells.Item(36, 4).Value = '1.175.31'
$ws.Cells.Item(36, 5).Value = '  +0.88%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '0.8916'
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -3.85%  '
$ws.Cells.Item(38, 5).Value = '  -5.07%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '0.5562'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -2.04%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '0.01555'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -2.42%  '
$ws.Cells.Item(41, 2).Value = 'PaxDollar'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -0.16%  '
$ws.Cells.Item(42, 2).Value = 'mCoin'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '2.548'
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -0.60%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '5.631'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -1.04%  '
$ws.Cells.Item(44, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '0.8046'
$c.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -4.87%  '
$ws.Cells.Item(45, 2).Value = 'Quant'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '99.45'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -1.08%  '
$ws.Cells.Item(46, 4).Value = '1.773.19'
$ws.Cells.Item(46, 5).Value = '  -1.36%  '
$ws.Cells.Item(47, 4).Value = '0.0₈109'
$ws.Cells.Item(47, 5).Value = '  -3.87%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '0.4513'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -0.49%  '
$ws.Cells.Item(49, 5).Value = '  -0.05%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '54.81'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -1.97%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '0.05040'
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -0.32%  '
